$d = $word.ActiveDocument

# 1. "-cse" -> "- cse"  (run has no special attrs; plain Find/Replace keeps structure intact)
$d.Content.Find.Execute("-cse", $true, $false, $false, $false, $false, $true, 1, $false, "- cse", 2)

# 2. Date string update. The target run carries w:rsidR="007C121E" and is immediately preceded
#    by a sibling run with identical rPr ("   "). A plain Find/Replace (or Range.Text=) on this
#    engine coalesces the edited run into that preceding sibling, which would corrupt the run
#    layout. Use Range.InsertXML with an explicit OOXML fragment so the run is rewritten in place
#    (including its original rsidR) without touching neighboring runs.
$marker = "  Date-2021-11-21 12:34:18.812834"
$p = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.IndexOf($marker) -ge 0) {
        $p = $cand
        break
    }
}
if ($p -eq $null) {
    throw "Paragraph containing the date stamp was not found"
}
$full = $p.Range.Text
$idx = $full.IndexOf($marker)
$rStart = $p.Range.Start + $idx
$rEnd = $rStart + $marker.Length
$dateRange = $d.Range($rStart, $rEnd)
$dateXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r w:rsidR="007C121E">' +
    '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '<w:t xml:space="preserve">  Date-2021-11-26 00:23:21.861267</w:t>' +
    '</w:r>' +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$dateRange.InsertXML($dateXml)

# 3. "unit assignment" -> "wertawetgfwa" (own run, no attrs)
$d.Content.Find.Execute("unit assignment", $true, $false, $false, $false, $false, $true, 1, $false, "wertawetgfwa", 2)

# 4. "submit it before 10 jan" -> "awertgawertg" (own run, no attrs)
$d.Content.Find.Execute("submit it before 10 jan", $true, $false, $false, $false, $false, $true, 1, $false, "awertgawertg", 2)

# 5. "teacher" -> "34234" (own run, no attrs)
$d.Content.Find.Execute("teacher", $true, $false, $false, $false, $false, $true, 1, $false, "34234", 2)
